{"js": "// The practice sheet has a single table; only 5 of its rows (0, 4, 9, 14, 19)\n// hold the multiplication expressions \u2014 the rest are blank rows left for\n// students to work in. Address each expression cell positionally (row/col)\n// rather than by searching for its old text, because several of the new\n// values happen to collide with old values used elsewhere in the table\n// (e.g. \"309\u00d74=\" is both a pre-existing value and a replacement value), so a\n// naive global find/replace would corrupt later cells.\nconst newValues = {\n  0: [\"245\u00d72=\", \"603\u00d72=\", \"222\u00d72=\", \"874\u00d79=\", \"760\u00d76=\"],\n  4: [\"309\u00d74=\", \"179\u00d76=\", \"353\u00d76=\", \"148\u00d73=\", \"805\u00d77=\"],\n  9: [\"878\u00d78=\", \"367\u00d74=\", \"314\u00d74=\", \"523\u00d72=\", \"350\u00d73=\"],\n  14: [\"678\u00d73=\", \"132\u00d73=\", \"608\u00d78=\", \"575\u00d79=\", \"634\u00d72=\"],\n  19: [\"528\u00d73=\", \"848\u00d79=\", \"382\u00d75=\", \"257\u00d76=\", \"985\u00d76=\"],\n};\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const [rowIndex, rowValues] of Object.entries(newValues)) {\n  const row = rows.items[Number(rowIndex)];\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n\n  for (let col = 0; col < rowValues.length; col++) {\n    const cell = cells.items[col];\n    const cellBody = cell.body;\n    const paragraphs = cellBody.paragraphs;\n    paragraphs.load(\"items\");\n    await context.sync();\n\n    const firstParagraph = paragraphs.items[0];\n    firstParagraph.load(\"text\");\n    await context.sync();\n\n    const range = firstParagraph.getRange();\n    range.insertText(rowValues[col], \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# The practice sheet has a single table; only 5 of its rows (1, 5, 10, 15, 20\n# in 1-based Word COM indexing) hold the multiplication expressions -- the\n# rest are blank rows left for students to work in. Address each expression\n# cell positionally (row/col) rather than via Find/Replace-by-old-text,\n# because several of the new values collide with old values used elsewhere\n# in the table (e.g. \"309x4=\" is both a pre-existing value and a replacement\n# value), so a naive global find/replace would corrupt later cells.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n    1  = @(\"245\u00d72=\", \"603\u00d72=\", \"222\u00d72=\", \"874\u00d79=\", \"760\u00d76=\")\n    5  = @(\"309\u00d74=\", \"179\u00d76=\", \"353\u00d76=\", \"148\u00d73=\", \"805\u00d77=\")\n    10 = @(\"878\u00d78=\", \"367\u00d74=\", \"314\u00d74=\", \"523\u00d72=\", \"350\u00d73=\")\n    15 = @(\"678\u00d73=\", \"132\u00d73=\", \"608\u00d78=\", \"575\u00d79=\", \"634\u00d72=\")\n    20 = @(\"528\u00d73=\", \"848\u00d79=\", \"382\u00d75=\", \"257\u00d76=\", \"985\u00d76=\")\n}\n\nforeach ($row in $newValues.Keys) {\n    $rowValues = $newValues[$row]\n    for ($col = 1; $col -le $rowValues.Length; $col++) {\n        $t.Cell($row, $col).Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
